$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $r2 = $ws.Range("$col" + "2").Value2
    $r3 = $ws.Range("$col" + "3").Value2
    $ws.Range("$col" + "2").Value2 = $r3
    $ws.Range("$col" + "3").Value2 = $r2
}

# Swap values between row 4 and row 5 for columns D, J, K, L, M, P
foreach ($col in $cols) {
    $r4 = $ws.Range("$col" + "4").Value2
    $r5 = $ws.Range("$col" + "5").Value2
    $ws.Range("$col" + "4").Value2 = $r5
    $ws.Range("$col" + "5").Value2 = $r4
}
